# Regenerate save_data column G ("K") values to use K instead of Strike#.
# This sets the recalculated K values for each row (rows 2-17) in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 4
    3  = 0
    4  = 5
    5  = 4
    6  = 1
    7  = 4
    8  = 2
    9  = 2
    10 = 3
    11 = 1
    12 = 0
    13 = 2
    14 = 3
    15 = 4
    16 = 1
    17 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
